# "(feature/Test): add method to get name school, detail account"
# The committed change re-assigns every class row (rows 2-29) on Sheet1
# to AcademicYearId = 11 (column D), and leaves the sheet's last-used
# selection on cell K24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D = AcademicYearId. Set it to 11 for every data row (2..29).
$ws.Range("D2:D29").Value = 11

# Move / record the active cell selection to K24 (was K5).
$ws.Range("K24").Select()
